$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Swap betting-odds data (columns F:V) between row pairs ---
$pairs = @(
    @(30, 31),
    @(42, 43),
    @(94, 95),
    @(97, 98),
    @(102, 103),
    @(128, 129),
    @(133, 134)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]
    $r1 = $ws.Range("F$row1`:V$row1")
    $r2 = $ws.Range("F$row2`:V$row2")
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value2 = $v2
    $r2.Value2 = $v1
}

# --- Step 2: Append 5 new match rows (138-142), copying formatting from row 137 ---
$ws.Range("A137:V137").Copy($ws.Range("A138:V138"))
$ws.Range("A137:V137").Copy($ws.Range("A139:V139"))
$ws.Range("A137:V137").Copy($ws.Range("A140:V140"))
$ws.Range("A137:V137").Copy($ws.Range("A141:V141"))
$ws.Range("A137:V137").Copy($ws.Range("A142:V142"))

$ws.Cells.Item(138, 1).Value = 137
$ws.Cells.Item(138, 2).Value = "portugal"
$ws.Cells.Item(138, 3).Value = "liga-portugal"
$ws.Cells.Item(138, 4).Value = "2023-2024"
$ws.Cells.Item(138, 5).Value = 45296.90625
$ws.Cells.Item(138, 6).Value = "Boavista"
$ws.Cells.Item(138, 7).Value = 1
$ws.Cells.Item(138, 8).Value = "FC Porto"
$ws.Cells.Item(138, 9).Value = 1
$ws.Cells.Item(138, 10).Value = 7.38
$ws.Cells.Item(138, 11).Value = "30/12/2024 18:13"
$ws.Cells.Item(138, 12).Value = 7.68
$ws.Cells.Item(138, 13).Value = "05/01/2024 21:44"
$ws.Cells.Item(138, 14).Value = 4.94
$ws.Cells.Item(138, 15).Value = "30/12/2024 18:13"
$ws.Cells.Item(138, 16).Value = 4.77
$ws.Cells.Item(138, 17).Value = "05/01/2024 21:44"
$ws.Cells.Item(138, 18).Value = 1.41
$ws.Cells.Item(138, 19).Value = "30/12/2024 18:13"
$ws.Cells.Item(138, 20).Value = 1.44
$ws.Cells.Item(138, 21).Value = "05/01/2024 21:42"
$ws.Cells.Item(138, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/boavista-fc-porto/xEwDuVVJ/"

$ws.Cells.Item(139, 1).Value = 138
$ws.Cells.Item(139, 2).Value = "portugal"
$ws.Cells.Item(139, 3).Value = "liga-portugal"
$ws.Cells.Item(139, 4).Value = "2023-2024"
$ws.Cells.Item(139, 5).Value = 45297.6875
$ws.Cells.Item(139, 6).Value = "SC Farense"
$ws.Cells.Item(139, 7).Value = 1
$ws.Cells.Item(139, 8).Value = "Gil Vicente"
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 2.16
$ws.Cells.Item(139, 11).Value = "30/12/2024 18:13"
$ws.Cells.Item(139, 12).Value = 2.3
$ws.Cells.Item(139, 13).Value = "06/01/2024 16:12"
$ws.Cells.Item(139, 14).Value = 3.65
$ws.Cells.Item(139, 15).Value = "30/12/2024 18:13"
$ws.Cells.Item(139, 16).Value = 3.34
$ws.Cells.Item(139, 17).Value = "06/01/2024 16:12"
$ws.Cells.Item(139, 18).Value = 3.32
$ws.Cells.Item(139, 19).Value = "30/12/2024 18:13"
$ws.Cells.Item(139, 20).Value = 3.41
$ws.Cells.Item(139, 21).Value = "06/01/2024 16:10"
$ws.Cells.Item(139, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/sc-farense-gil-vicente/YqPrz9gs/"

$ws.Cells.Item(140, 1).Value = 139
$ws.Cells.Item(140, 2).Value = "portugal"
$ws.Cells.Item(140, 3).Value = "liga-portugal"
$ws.Cells.Item(140, 4).Value = "2023-2024"
$ws.Cells.Item(140, 5).Value = 45297.6875
$ws.Cells.Item(140, 6).Value = "Estrela"
$ws.Cells.Item(140, 7).Value = 1
$ws.Cells.Item(140, 8).Value = "Vizela"
$ws.Cells.Item(140, 9).Value = 1
$ws.Cells.Item(140, 10).Value = 2.45
$ws.Cells.Item(140, 11).Value = "30/12/2024 18:13"
$ws.Cells.Item(140, 12).Value = 2.93
$ws.Cells.Item(140, 13).Value = "06/01/2024 16:29"
$ws.Cells.Item(140, 14).Value = 3.27
$ws.Cells.Item(140, 15).Value = "30/12/2024 18:13"
$ws.Cells.Item(140, 16).Value = 3.3
$ws.Cells.Item(140, 17).Value = "06/01/2024 16:29"
$ws.Cells.Item(140, 18).Value = 3.07
$ws.Cells.Item(140, 19).Value = "30/12/2024 18:13"
$ws.Cells.Item(140, 20).Value = 2.61
$ws.Cells.Item(140, 21).Value = "06/01/2024 16:29"
$ws.Cells.Item(140, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/estrela-da-amadora-vizela/G4xHvkoQ/"

$ws.Cells.Item(141, 1).Value = 140
$ws.Cells.Item(141, 2).Value = "portugal"
$ws.Cells.Item(141, 3).Value = "liga-portugal"
$ws.Cells.Item(141, 4).Value = "2023-2024"
$ws.Cells.Item(141, 5).Value = 45297.79166666666
$ws.Cells.Item(141, 6).Value = "Arouca"
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = "Benfica"
$ws.Cells.Item(141, 9).Value = 3
$ws.Cells.Item(141, 10).Value = 6.29
$ws.Cells.Item(141, 11).Value = "30/12/2024 19:12"
$ws.Cells.Item(141, 12).Value = 7.24
$ws.Cells.Item(141, 13).Value = "06/01/2024 18:59"
$ws.Cells.Item(141, 14).Value = 4.94
$ws.Cells.Item(141, 15).Value = "30/12/2024 19:12"
$ws.Cells.Item(141, 16).Value = 4.66
$ws.Cells.Item(141, 17).Value = "06/01/2024 18:59"
$ws.Cells.Item(141, 18).Value = 1.46
$ws.Cells.Item(141, 19).Value = "30/12/2024 19:12"
$ws.Cells.Item(141, 20).Value = 1.47
$ws.Cells.Item(141, 21).Value = "06/01/2024 18:57"
$ws.Cells.Item(141, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/arouca-benfica/0dfUGogJ/"

$ws.Cells.Item(142, 1).Value = 141
$ws.Cells.Item(142, 2).Value = "portugal"
$ws.Cells.Item(142, 3).Value = "liga-portugal"
$ws.Cells.Item(142, 4).Value = "2023-2024"
$ws.Cells.Item(142, 5).Value = 45297.89583333334
$ws.Cells.Item(142, 6).Value = "Braga"
$ws.Cells.Item(142, 7).Value = 1
$ws.Cells.Item(142, 8).Value = "Vitoria Guimaraes"
$ws.Cells.Item(142, 9).Value = 1
$ws.Cells.Item(142, 10).Value = 1.69
$ws.Cells.Item(142, 11).Value = "30/12/2024 21:42"
$ws.Cells.Item(142, 12).Value = 1.63
$ws.Cells.Item(142, 13).Value = "06/01/2024 21:23"
$ws.Cells.Item(142, 14).Value = 4.06
$ws.Cells.Item(142, 15).Value = "30/12/2024 21:42"
$ws.Cells.Item(142, 16).Value = 4.29
$ws.Cells.Item(142, 17).Value = "06/01/2024 21:25"
$ws.Cells.Item(142, 18).Value = 4.83
$ws.Cells.Item(142, 19).Value = "30/12/2024 21:42"
$ws.Cells.Item(142, 20).Value = 5.48
$ws.Cells.Item(142, 21).Value = "06/01/2024 21:25"
$ws.Cells.Item(142, 22).Value = "https://www.betexplorer.com/football/portugal/liga-portugal/braga-vitoria-guimaraes/M5eQHRwD/"
